$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B3 value
$ws.Range("B3").Value = 864831.5978229975

# Clear D3 (becomes an empty inline string cell)
$ws.Range("D3").Value = ""

# Update C4 value
$ws.Range("C4").Value = 66.16930873157268

# Update C5 value
$ws.Range("C5").Value = 1864.133311530407

# Row 7 "Other" becomes "Biogas" with updated D7 value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 801.5679477565395

# New row 8 "Other" with D8 value, copying style from row 7's A7 cell
$ws.Range("A8").Value = "Other"
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D8").Value = 877.9100423599186
